{"js": "// Remove the trailing \"Ver no Jupiter...\" line, the copyright/footer line,\n// and the blank paragraph that separates them from the bibliography text,\n// mirroring the upstream Jekyll site rebuild that dropped this footer block.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// Find the index of the \"Ver no Jupiter...\" paragraph; the blank paragraph\n// immediately preceding it (the separator after the bibliography) and the\n// copyright paragraph immediately following it are removed together.\nlet jupiterIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === jupiterText) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex !== -1) {\n  const toDelete = [];\n  // Preceding blank paragraph (separator), if present.\n  if (jupiterIndex - 1 >= 0 && paragraphs.items[jupiterIndex - 1].text === \"\") {\n    toDelete.push(paragraphs.items[jupiterIndex - 1]);\n  }\n  toDelete.push(paragraphs.items[jupiterIndex]);\n  if (\n    jupiterIndex + 1 < paragraphs.items.length &&\n    paragraphs.items[jupiterIndex + 1].text === copyrightText\n  ) {\n    toDelete.push(paragraphs.items[jupiterIndex + 1]);\n  }\n\n  for (const para of toDelete) {\n    para.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" line, the copyright/footer line,\n# and the blank paragraph that separates them from the bibliography text,\n# mirroring the upstream Jekyll site rebuild that dropped this footer block.\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightMarker = \"Contact: luizeleno@usp.br\"\n\n$jupiterIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.TrimEnd([char]13, [char]7) -eq $jupiterText) {\n        $jupiterIndex = $i\n        break\n    }\n}\n\nif ($jupiterIndex -ne -1) {\n    # Delete from the end backward so earlier indices stay valid.\n    $copyrightIndex = $jupiterIndex + 1\n    if ($copyrightIndex -le $count) {\n        $copyrightText = $d.Paragraphs.Item($copyrightIndex).Range.Text\n        if ($copyrightText -like \"*$copyrightMarker*\") {\n            $d.Paragraphs.Item($copyrightIndex).Range.Delete()\n        }\n    }\n\n    $d.Paragraphs.Item($jupiterIndex).Range.Delete()\n\n    $blankIndex = $jupiterIndex - 1\n    if ($blankIndex -ge 1) {\n        $blankText = $d.Paragraphs.Item($blankIndex).Range.Text\n        if ($blankText.TrimEnd([char]13, [char]7) -eq \"\") {\n            $d.Paragraphs.Item($blankIndex).Range.Delete()\n        }\n    }\n}\n"}
